$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first two fixtures (Man Utd v Aston Villa, MU Women v Birmingham WFC)
$ws.Rows("1:2").Delete()

# Insert a new fixture row before "Manchester United v Southampton" (now row 3)
# so it sits right after "MU Women v Tottenham Hotspur Women"
$ws.Rows("3:3").Insert()
$ws.Range("A3").Value = "Manchester United v Middlesbrough (FA Cup) "
$ws.Range("B3").Value = "04 FebFri20:00"
